# Daily attendance processing - 2026-01-02 06:05:07
# Swap the order of "Recorded By" values from "System, dnasr281@gmail.com"
# to "dnasr281@gmail.com, System" for the affected session rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = @(
    "G8","G9","G10","G12","G14","G15","G17","G18",
    "G34","G35","G36","G38","G40","G41","G43","G44",
    "G60","G61","G62","G64","G66","G67","G69","G70",
    "G86","G87","G88","G90","G92","G93","G95","G96",
    "G112","G113","G114","G116","G118","G119","G121","G122",
    "G138","G139","G140","G142","G144","G145","G147","G148",
    "G164","G167","G170","G174",
    "G191","G194","G197","G201",
    "G218","G221","G224","G228",
    "G245","G248","G251","G255",
    "G272","G275","G278","G282",
    "G299","G302","G305","G309"
)

foreach ($addr in $cells) {
    $ws.Range($addr).Value = "dnasr281@gmail.com, System"
}
